# "Logged Week 15 and simulated Week 16"
#
# Week 15 (RB sheet) is now fully logged: a new player row (G.Nabers) is
# appended with zeroed stat columns, and the RB sheet becomes the active
# sheet/tab (simulation has moved on to Week 16, the next sheet in the
# workflow), leaving the WR sheet no longer selected.

$wb = $excel.ActiveWorkbook

# --- RB sheet: add the new player row for the week just logged ---
$rb = $wb.Worksheets.Item("RB")

$rb.Range("A6").Value = "G.Nabers"
$rb.Range("B6:J6").Value = 0

# --- Make RB the active sheet/tab, with the selection parked just past
#     the new data (mirrors the "ready for next week" cursor position) ---
$rb.Activate()
$rb.Range("J7").Select()
